$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 35.88321066666667
$ws.Cells.Item(2, 8).Value = 107.649632
$ws.Cells.Item(2, 9).Value = 0.08317795499144418
$ws.Cells.Item(2, 10).Value = 0.08448843719082051
$ws.Cells.Item(2, 13).Value = 19.77408333333333
$ws.Cells.Item(2, 14).Value = 59.32225
$ws.Cells.Item(2, 15).Value = 0.3380388258879848
$ws.Cells.Item(2, 16).Value = 0.339186328349942
$ws.Cells.Item(2, 17).Value = 709.5575979902223
$ws.Cells.Item(2, 18).Value = 6386.018381912001
$ws.Cells.Item(2, 19).Value = 0.02811737824507143
$ws.Cells.Item(2, 20).Value = 0.0286573227987791
$ws.Cells.Item(3, 7).Value = 35.88321066666667
$ws.Cells.Item(3, 8).Value = 107.649632
$ws.Cells.Item(3, 9).Value = 0.08317795499144418
$ws.Cells.Item(3, 10).Value = 0.08448843719082051
$ws.Cells.Item(3, 15).Value = 0.3069959581674471
$ws.Cells.Item(3, 16).Value = 0.3080380828904952
$ws.Cells.Item(3, 17).Value = 644.3973235849032
$ws.Cells.Item(3, 18).Value = 5799.575912264129
$ws.Cells.Item(3, 19).Value = 0.02553529599100719
$ws.Cells.Item(3, 20).Value = 0.02602565621867437
$ws.Cells.Item(4, 7).Value = 35.88321066666667
$ws.Cells.Item(4, 8).Value = 107.649632
$ws.Cells.Item(4, 9).Value = 0.08317795499144418
$ws.Cells.Item(4, 10).Value = 0.08448843719082051
$ws.Cells.Item(4, 13).Value = 12.46730333333333
$ws.Cells.Item(4, 14).Value = 37.40191
$ws.Cells.Item(4, 15).Value = 0.2131290998296268
$ws.Cells.Item(4, 16).Value = 0.2138525852639604
$ws.Cells.Item(4, 17).Value = 447.3668719552356
$ws.Cells.Item(4, 18).Value = 4026.30184759712
$ws.Cells.Item(4, 19).Value = 0.01772764267299571
$ws.Cells.Item(4, 20).Value = 0.01806807071816871
$ws.Cells.Item(5, 7).Value = 35.88321066666667
$ws.Cells.Item(5, 8).Value = 107.649632
$ws.Cells.Item(5, 9).Value = 0.08317795499144418
$ws.Cells.Item(5, 10).Value = 0.08448843719082051
$ws.Cells.Item(5, 13).Value = 0.593699
$ws.Cells.Item(5, 14).Value = 1.187398
$ws.Cells.Item(5, 15).Value = 0.01014931056513554
$ws.Cells.Item(5, 16).Value = 0.006789175527058808
$ws.Cells.Item(5, 17).Value = 21.30382628958933
$ws.Cells.Item(5, 18).Value = 127.822957737536
$ws.Cells.Item(5, 19).Value = 0.0008441988973810331
$ws.Cells.Item(5, 20).Value = 0.0005736068300953638
$ws.Cells.Item(6, 7).Value = 35.88321066666667
$ws.Cells.Item(6, 8).Value = 107.649632
$ws.Cells.Item(6, 9).Value = 0.08317795499144418
$ws.Cells.Item(6, 10).Value = 0.08448843719082051
$ws.Cells.Item(6, 13).Value = 7.703215333333333
$ws.Cells.Item(6, 14).Value = 23.109646
$ws.Cells.Item(6, 15).Value = 0.1316868055498057
$ws.Cells.Item(6, 16).Value = 0.1321338279685434
$ws.Cells.Item(6, 17).Value = 276.4160986166969
$ws.Cells.Item(6, 18).Value = 2487.744887550272
$ws.Cells.Item(6, 19).Value = 0.0109534391849888
$ws.Cells.Item(6, 20).Value = 0.01116378062510296
$ws.Cells.Item(7, 9).Value = 0.03522729558434242
$ws.Cells.Item(7, 10).Value = 0.03578230735158529
$ws.Cells.Item(7, 13).Value = 19.77408333333333
$ws.Cells.Item(7, 14).Value = 59.32225
$ws.Cells.Item(7, 15).Value = 0.3380388258879848
$ws.Cells.Item(7, 16).Value = 0.339186328349942
$ws.Cells.Item(7, 17).Value = 300.5098555391111
$ws.Cells.Item(7, 18).Value = 2704.588699852
$ws.Cells.Item(7, 19).Value = 0.0119081936385401
$ws.Cells.Item(7, 20).Value = 0.01213686945047335
$ws.Cells.Item(8, 9).Value = 0.03522729558434242
$ws.Cells.Item(8, 10).Value = 0.03578230735158529
$ws.Cells.Item(8, 15).Value = 0.3069959581674471
$ws.Cells.Item(8, 16).Value = 0.3080380828904952
$ws.Cells.Item(8, 19).Value = 0.01081463736156308
$ws.Cells.Item(8, 20).Value = 0.01102231335798081
$ws.Cells.Item(9, 9).Value = 0.03522729558434242
$ws.Cells.Item(9, 10).Value = 0.03578230735158529
$ws.Cells.Item(9, 13).Value = 12.46730333333333
$ws.Cells.Item(9, 14).Value = 37.40191
$ws.Cells.Item(9, 15).Value = 0.2131290998296268
$ws.Cells.Item(9, 16).Value = 0.2138525852639604
$ws.Cells.Item(9, 17).Value = 189.4675702790578
$ws.Cells.Item(9, 18).Value = 1705.20813251152
$ws.Cells.Item(9, 19).Value = 0.007507961797323087
$ws.Cells.Item(9, 20).Value = 0.007652138933846133
$ws.Cells.Item(10, 9).Value = 0.03522729558434242
$ws.Cells.Item(10, 10).Value = 0.03578230735158529
$ws.Cells.Item(10, 13).Value = 0.593699
$ws.Cells.Item(10, 14).Value = 1.187398
$ws.Cells.Item(10, 15).Value = 0.01014931056513554
$ws.Cells.Item(10, 16).Value = 0.006789175527058808
$ws.Cells.Item(10, 17).Value = 9.022537111642665
$ws.Cells.Item(10, 18).Value = 54.13522266985599
$ws.Cells.Item(10, 19).Value = 0.0003575327632553192
$ws.Cells.Item(10, 20).Value = 0.0002429323653730794
$ws.Cells.Item(11, 9).Value = 0.03522729558434242
$ws.Cells.Item(11, 10).Value = 0.03578230735158529
$ws.Cells.Item(11, 13).Value = 7.703215333333333
$ws.Cells.Item(11, 14).Value = 23.109646
$ws.Cells.Item(11, 15).Value = 0.1316868055498057
$ws.Cells.Item(11, 16).Value = 0.1321338279685434
$ws.Cells.Item(11, 17).Value = 117.0669753932124
$ws.Cells.Item(11, 18).Value = 1053.602778538912
$ws.Cells.Item(11, 19).Value = 0.004638970023660831
$ws.Cells.Item(11, 20).Value = 0.004728053243911916
$ws.Cells.Item(12, 7).Value = 177.70077
$ws.Cells.Item(12, 8).Value = 533.10231
$ws.Cells.Item(12, 9).Value = 0.4119137160358794
$ws.Cells.Item(12, 10).Value = 0.4184034835782469
$ws.Cells.Item(12, 13).Value = 19.77408333333333
$ws.Cells.Item(12, 14).Value = 59.32225
$ws.Cells.Item(12, 15).Value = 0.3380388258879848
$ws.Cells.Item(12, 16).Value = 0.339186328349942
$ws.Cells.Item(12, 17).Value = 3513.8698343775
$ws.Cells.Item(12, 18).Value = 31624.8285093975
$ws.Cells.Item(12, 19).Value = 0.1392428289359254
$ws.Cells.Item(12, 20).Value = 0.1419167413637308
$ws.Cells.Item(13, 7).Value = 177.70077
$ws.Cells.Item(13, 8).Value = 533.10231
$ws.Cells.Item(13, 9).Value = 0.4119137160358794
$ws.Cells.Item(13, 10).Value = 0.4184034835782469
$ws.Cells.Item(13, 15).Value = 0.3069959581674471
$ws.Cells.Item(13, 16).Value = 0.3080380828904952
$ws.Cells.Item(13, 17).Value = 3191.183243068861
$ws.Cells.Item(13, 18).Value = 28720.64918761974
$ws.Cells.Item(13, 19).Value = 0.1264558459367485
$ws.Cells.Item(13, 20).Value = 0.128884206956148
$ws.Cells.Item(14, 7).Value = 177.70077
$ws.Cells.Item(14, 8).Value = 533.10231
$ws.Cells.Item(14, 9).Value = 0.4119137160358794
$ws.Cells.Item(14, 10).Value = 0.4184034835782469
$ws.Cells.Item(14, 13).Value = 12.46730333333333
$ws.Cells.Item(14, 14).Value = 37.40191
$ws.Cells.Item(14, 15).Value = 0.2131290998296268
$ws.Cells.Item(14, 16).Value = 0.2138525852639604
$ws.Cells.Item(14, 17).Value = 2215.4494021569
$ws.Cells.Item(14, 18).Value = 19939.0446194121
$ws.Cells.Item(14, 19).Value = 0.08779079950620348
$ws.Cells.Item(14, 20).Value = 0.08947666664665511
$ws.Cells.Item(15, 7).Value = 177.70077
$ws.Cells.Item(15, 8).Value = 533.10231
$ws.Cells.Item(15, 9).Value = 0.4119137160358794
$ws.Cells.Item(15, 10).Value = 0.4184034835782469
$ws.Cells.Item(15, 13).Value = 0.593699
$ws.Cells.Item(15, 14).Value = 1.187398
$ws.Cells.Item(15, 15).Value = 0.01014931056513554
$ws.Cells.Item(15, 16).Value = 0.006789175527058808
$ws.Cells.Item(15, 17).Value = 105.50076944823
$ws.Cells.Item(15, 18).Value = 633.00461668938
$ws.Cells.Item(15, 19).Value = 0.004180640230087193
$ws.Cells.Item(15, 20).Value = 0.002840614691145585
$ws.Cells.Item(16, 7).Value = 177.70077
$ws.Cells.Item(16, 8).Value = 533.10231
$ws.Cells.Item(16, 9).Value = 0.4119137160358794
$ws.Cells.Item(16, 10).Value = 0.4184034835782469
$ws.Cells.Item(16, 13).Value = 7.703215333333333
$ws.Cells.Item(16, 14).Value = 23.109646
$ws.Cells.Item(16, 15).Value = 0.1316868055498057
$ws.Cells.Item(16, 16).Value = 0.1321338279685434
$ws.Cells.Item(16, 17).Value = 1368.86729620914
$ws.Cells.Item(16, 18).Value = 12319.80566588226
$ws.Cells.Item(16, 19).Value = 0.05424360142691475
$ws.Cells.Item(16, 20).Value = 0.05528525392056733
$ws.Cells.Item(17, 7).Value = 20.074196
$ws.Cells.Item(17, 8).Value = 40.148392
$ws.Cells.Item(17, 9).Value = 0.04653236263856699
$ws.Cells.Item(17, 10).Value = 0.0315103250497358
$ws.Cells.Item(17, 13).Value = 19.77408333333333
$ws.Cells.Item(17, 14).Value = 59.32225
$ws.Cells.Item(17, 15).Value = 0.3380388258879848
$ws.Cells.Item(17, 16).Value = 0.339186328349942
$ws.Cells.Item(17, 17).Value = 396.9488245536667
$ws.Cells.Item(17, 18).Value = 2381.692947322
$ws.Cells.Item(17, 19).Value = 0.01572974523213512
$ws.Cells.Item(17, 20).Value = 0.01068787145873309
$ws.Cells.Item(18, 7).Value = 20.074196
$ws.Cells.Item(18, 8).Value = 40.148392
$ws.Cells.Item(18, 9).Value = 0.04653236263856699
$ws.Cells.Item(18, 10).Value = 0.0315103250497358
$ws.Cells.Item(18, 15).Value = 0.3069959581674471
$ws.Cells.Item(18, 16).Value = 0.3080380828904952
$ws.Cells.Item(18, 17).Value = 360.4961188028614
$ws.Cells.Item(18, 18).Value = 2162.976712817168
$ws.Cells.Item(18, 19).Value = 0.01428524725402199
$ws.Cells.Item(18, 20).Value = 0.009706380119576963
$ws.Cells.Item(19, 7).Value = 20.074196
$ws.Cells.Item(19, 8).Value = 40.148392
$ws.Cells.Item(19, 9).Value = 0.04653236263856699
$ws.Cells.Item(19, 10).Value = 0.0315103250497358
$ws.Cells.Item(19, 13).Value = 12.46730333333333
$ws.Cells.Item(19, 14).Value = 37.40191
$ws.Cells.Item(19, 15).Value = 0.2131290998296268
$ws.Cells.Item(19, 16).Value = 0.2138525852639604
$ws.Cells.Item(19, 17).Value = 250.2710907047867
$ws.Cells.Item(19, 18).Value = 1501.62654422872
$ws.Cells.Item(19, 19).Value = 0.009917400562103541
$ws.Cells.Item(19, 20).Value = 0.006738564474393733
$ws.Cells.Item(20, 7).Value = 20.074196
$ws.Cells.Item(20, 8).Value = 40.148392
$ws.Cells.Item(20, 9).Value = 0.04653236263856699
$ws.Cells.Item(20, 10).Value = 0.0315103250497358
$ws.Cells.Item(20, 13).Value = 0.593699
$ws.Cells.Item(20, 14).Value = 1.187398
$ws.Cells.Item(20, 15).Value = 0.01014931056513554
$ws.Cells.Item(20, 16).Value = 0.006789175527058808
$ws.Cells.Item(20, 17).Value = 11.918030091004
$ws.Cells.Item(20, 18).Value = 47.672120364016
$ws.Cells.Item(20, 19).Value = 0.0004722713997483264
$ws.Cells.Item(20, 20).Value = 0.0002139291276773344
$ws.Cells.Item(21, 7).Value = 20.074196
$ws.Cells.Item(21, 8).Value = 40.148392
$ws.Cells.Item(21, 9).Value = 0.04653236263856699
$ws.Cells.Item(21, 10).Value = 0.0315103250497358
$ws.Cells.Item(21, 13).Value = 7.703215333333333
$ws.Cells.Item(21, 14).Value = 23.109646
$ws.Cells.Item(21, 15).Value = 0.1316868055498057
$ws.Cells.Item(21, 16).Value = 0.1321338279685434
$ws.Cells.Item(21, 17).Value = 154.6358544315387
$ws.Cells.Item(21, 18).Value = 927.8151265892319
$ws.Cells.Item(21, 19).Value = 0.006127698190558017
$ws.Cells.Item(21, 20).Value = 0.004163579869354672
$ws.Cells.Item(22, 7).Value = 182.547562
$ws.Cells.Item(22, 8).Value = 547.642686
$ws.Cells.Item(22, 9).Value = 0.423148670749767
$ws.Cells.Item(22, 10).Value = 0.4298154468296114
$ws.Cells.Item(22, 13).Value = 19.77408333333333
$ws.Cells.Item(22, 14).Value = 59.32225
$ws.Cells.Item(22, 15).Value = 0.3380388258879848
$ws.Cells.Item(22, 16).Value = 0.339186328349942
$ws.Cells.Item(22, 17).Value = 3609.710703284833
$ws.Cells.Item(22, 18).Value = 32487.3963295635
$ws.Cells.Item(22, 19).Value = 0.1430406798363127
$ws.Cells.Item(22, 20).Value = 0.1457875232782256
$ws.Cells.Item(23, 7).Value = 182.547562
$ws.Cells.Item(23, 8).Value = 547.642686
$ws.Cells.Item(23, 9).Value = 0.423148670749767
$ws.Cells.Item(23, 10).Value = 0.4298154468296114
$ws.Cells.Item(23, 15).Value = 0.3069959581674471
$ws.Cells.Item(23, 16).Value = 0.3080380828904952
$ws.Cells.Item(23, 17).Value = 3278.222828845783
$ws.Cells.Item(23, 18).Value = 29504.00545961205
$ws.Cells.Item(23, 19).Value = 0.1299049316241063
$ws.Cells.Item(23, 20).Value = 0.1323995262381151
$ws.Cells.Item(24, 7).Value = 182.547562
$ws.Cells.Item(24, 8).Value = 547.642686
$ws.Cells.Item(24, 9).Value = 0.423148670749767
$ws.Cells.Item(24, 10).Value = 0.4298154468296114
$ws.Cells.Item(24, 13).Value = 12.46730333333333
$ws.Cells.Item(24, 14).Value = 37.40191
$ws.Cells.Item(24, 15).Value = 0.2131290998296268
$ws.Cells.Item(24, 16).Value = 0.2138525852639604
$ws.Cells.Item(24, 17).Value = 2275.875828214474
$ws.Cells.Item(24, 18).Value = 20482.88245393026
$ws.Cells.Item(24, 19).Value = 0.09018529529100099
$ws.Cells.Item(24, 20).Value = 0.09191714449089673
$ws.Cells.Item(25, 7).Value = 182.547562
$ws.Cells.Item(25, 8).Value = 547.642686
$ws.Cells.Item(25, 9).Value = 0.423148670749767
$ws.Cells.Item(25, 10).Value = 0.4298154468296114
$ws.Cells.Item(25, 13).Value = 0.593699
$ws.Cells.Item(25, 14).Value = 1.187398
$ws.Cells.Item(25, 15).Value = 0.01014931056513554
$ws.Cells.Item(25, 16).Value = 0.006789175527058808
$ws.Cells.Item(25, 17).Value = 108.378305011838
$ws.Cells.Item(25, 18).Value = 650.269830071028
$ws.Cells.Item(25, 19).Value = 0.004294667274663672
$ws.Cells.Item(25, 20).Value = 0.002918092512767444
$ws.Cells.Item(26, 7).Value = 182.547562
$ws.Cells.Item(26, 8).Value = 547.642686
$ws.Cells.Item(26, 9).Value = 0.423148670749767
$ws.Cells.Item(26, 10).Value = 0.4298154468296114
$ws.Cells.Item(26, 13).Value = 7.703215333333333
$ws.Cells.Item(26, 14).Value = 23.109646
$ws.Cells.Item(26, 15).Value = 0.1316868055498057
$ws.Cells.Item(26, 16).Value = 0.1321338279685434
$ws.Cells.Item(26, 17).Value = 1406.203178661017
$ws.Cells.Item(26, 18).Value = 12655.82860794916
$ws.Cells.Item(26, 19).Value = 0.05572309672368334
$ws.Cells.Item(26, 20).Value = 0.05679316030960647
